# Regenerate orders with updated distance/sizes:
#   D64 -> D69, D80 -> D86, D51 -> D55, S30 -> S31
# Applied as a text substitution across every string-valued cell
# in the worksheet's used range (these codes are embedded inside
# larger tokens such as "Face12_D64_S30" or "Face12_D64_S30_l.png",
# as well as appearing standalone in the Distance/Size columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string]) {
            if ($v.Contains("D64") -or $v.Contains("D80") -or $v.Contains("D51") -or $v.Contains("S30")) {
                $newV = $v.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
                $cell.Value = $newV
            }
        }
    }
}
